$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '69.967.48'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  -1.50%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.748.69'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +2.30%  '

$ws.Range("E4").Value = '  -0.06%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '622.61'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +0.49%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '179.96'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -0.89%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.747.48'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +2.43%  '

$ws.Range("E8").Value = '  +0.00%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.535'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -1.47%  '

$ws.Range("E10").Value = '  +3.07%  '

$ws.Range("E11").Value = '  -5.40%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.487'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -3.25%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '41.15'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +1.79%  '

$ws.Range("E14").Value = '  +2.05%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.361.54'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +1.99%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.747.01'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +2.43%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '69.992.33'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -1.46%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.122'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -1.42%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.59'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +0.55%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '16.79'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -0.55%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '506.78'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -2.77%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.41'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +2.16%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.725'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -2.33%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.54'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +0.13%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '86.87'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -1.94%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '13.13'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -2.80%  '

$ws.Range("E27").Value = '  +0.56%  '

$ws.Range("E28").Value = '  +24.70%  '

$ws.Range("E29").Value = '  -0.02%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.48'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -2.82%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.93'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +0.26%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.92'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -3.41%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '31.31'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -0.93%  '

$ws.Range("E34").Value = '  -0.24%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.999'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -0.06%  '

$ws.Range("E36").Value = '  +3.23%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.21'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +1.28%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.335'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -4.23%  '

$ws.Range("E39").Value = '  +1.06%  '

$ws.Range("E40").Value = '  -4.43%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '50.40'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -1.62%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '45.04'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -1.28%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '424.87'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -1.51%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.72'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -1.28%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.84'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +0.96%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.997.22'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -3.75%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0364'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -1.63%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '27.30'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -3.98%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '138.20'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -1.33%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.49'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +0.74%  '
